$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

$ws_ALC.Range("H40").Value = 3067.52
$ws_ALC.Range("I40").Value = 1606.2858
$ws_ALC.Range("J40").Value = 4927.273
$ws_ALC.Range("K40").Value = 1606.2858
$ws_ALC.Range("L40").Value = 4927.273
$ws_ALC.Range("M40").Value = -1431.2858
$ws_ALC.Range("N40").Value = -5277.273

$ws_ALC.Range("H62").Value = 3700.7144
$ws_ALC.Range("I62").Value = 3781
$ws_ALC.Range("J62").Value = 3500
$ws_ALC.Range("K62").Value = 3781
$ws_ALC.Range("L62").Value = 3500
$ws_ALC.Range("M62").Value = -3157
$ws_ALC.Range("N62").Value = -4748

$ws_ALC.Range("H65").Value = 3700.7144
$ws_ALC.Range("I65").Value = 3781
$ws_ALC.Range("J65").Value = 3500
$ws_ALC.Range("K65").Value = 18905
$ws_ALC.Range("L65").Value = 17500
$ws_ALC.Range("M65").Value = -15785
$ws_ALC.Range("N65").Value = -23740

$ws_ALC.Range("H129").Value = 872.4483
$ws_ALC.Range("I129").Value = 213.16667
$ws_ALC.Range("J129").Value = 1044.4348
$ws_ALC.Range("K129").Value = 639.50001
$ws_ALC.Range("L129").Value = 3133.3044
$ws_ALC.Range("M129").Value = 4360.49999
$ws_ALC.Range("N129").Value = -13133.3044

$ws_ALC.Range("H132").Value = 1570.5714
$ws_ALC.Range("I132").Value = 848.6667
$ws_ALC.Range("K132").Value = 2546.0001
$ws_ALC.Range("M132").Value = -16.0001000000002

$ws_ALC.Range("H138").Value = 4351929.5
$ws_ALC.Range("I138").Value = 2501.353
$ws_ALC.Range("J138").Value = 6901594
$ws_ALC.Range("K138").Value = 7504.059
$ws_ALC.Range("L138").Value = 20704782
$ws_ALC.Range("M138").Value = -2364.059
$ws_ALC.Range("N138").Value = -20715062

$ws_ARM.Range("H122").Value = 1427.3334
$ws_ARM.Range("I122").Value = 1149.037
$ws_ARM.Range("J122").Value = 2053.5
$ws_ARM.Range("K122").Value = 3447.111
$ws_ARM.Range("L122").Value = 6160.5
$ws_ARM.Range("M122").Value = -997.1109999999999
$ws_ARM.Range("N122").Value = -11060.5

$ws_ARM.Range("H132").Value = 1685.5227
$ws_ARM.Range("I132").Value = 1414.5667
$ws_ARM.Range("K132").Value = 4243.7001
$ws_ARM.Range("M132").Value = -1713.7001

$ws_BSM.Range("H20").Value = 2102.25
$ws_BSM.Range("I20").Value = 1000
$ws_BSM.Range("J20").Value = 2469.6667
$ws_BSM.Range("K20").Value = 1000
$ws_BSM.Range("L20").Value = 2469.6667
$ws_BSM.Range("M20").Value = -753
$ws_BSM.Range("N20").Value = -2963.6667

$ws_BSM.Range("H86").Value = 1782.2413
$ws_BSM.Range("I86").Value = 1685.6818
$ws_BSM.Range("J86").Value = 2085.7144
$ws_BSM.Range("K86").Value = 1685.6818
$ws_BSM.Range("L86").Value = 2085.7144
$ws_BSM.Range("M86").Value = -562.6818000000001
$ws_BSM.Range("N86").Value = -4331.7144

$ws_BSM.Range("H89").Value = 1782.2413
$ws_BSM.Range("I89").Value = 1685.6818
$ws_BSM.Range("J89").Value = 2085.7144
$ws_BSM.Range("K89").Value = 8428.409
$ws_BSM.Range("L89").Value = 10428.572
$ws_BSM.Range("M89").Value = -2812.409
$ws_BSM.Range("N89").Value = -21660.572

$ws_BSM.Range("H134").Value = 4289.246
$ws_BSM.Range("I134").Value = 4472.7803
$ws_BSM.Range("J134").Value = 3975.7083
$ws_BSM.Range("K134").Value = 13418.3409
$ws_BSM.Range("L134").Value = 11927.1249
$ws_BSM.Range("M134").Value = -10883.3409
$ws_BSM.Range("N134").Value = -16997.1249

$ws_BSM.Range("H138").Value = 36056.5
$ws_BSM.Range("J138").Value = 36056.5
$ws_BSM.Range("L138").Value = 36056.5
$ws_BSM.Range("N138").Value = -46336.5

$ws_CRP.Range("H31").Value = 11629471
$ws_CRP.Range("I31").Value = 1158.4375
$ws_CRP.Range("J31").Value = 18520322
$ws_CRP.Range("K31").Value = 1158.4375
$ws_CRP.Range("L31").Value = 18520322
$ws_CRP.Range("M31").Value = -863.4375
$ws_CRP.Range("N31").Value = -18520912

$ws_CRP.Range("H34").Value = 11629471
$ws_CRP.Range("I34").Value = 1158.4375
$ws_CRP.Range("J34").Value = 18520322
$ws_CRP.Range("K34").Value = 1158.4375
$ws_CRP.Range("L34").Value = 18520322
$ws_CRP.Range("M34").Value = -956.4375
$ws_CRP.Range("N34").Value = -18520726

$ws_CRP.Range("H69").Value = 9944.571
$ws_CRP.Range("I69").Value = 1922.6
$ws_CRP.Range("J69").Value = 29999.5
$ws_CRP.Range("K69").Value = 1922.6
$ws_CRP.Range("L69").Value = 29999.5
$ws_CRP.Range("M69").Value = -1173.6
$ws_CRP.Range("N69").Value = -31497.5

$ws_CRP.Range("H72").Value = 9944.571
$ws_CRP.Range("I72").Value = 1922.6
$ws_CRP.Range("J72").Value = 29999.5
$ws_CRP.Range("K72").Value = 5767.799999999999
$ws_CRP.Range("L72").Value = 89998.5
$ws_CRP.Range("M72").Value = -2023.799999999999
$ws_CRP.Range("N72").Value = -97486.5

$ws_CRP.Range("H86").Value = 14886
$ws_CRP.Range("I86").Value = 4501.75
$ws_CRP.Range("K86").Value = 4501.75
$ws_CRP.Range("M86").Value = -3378.75

$ws_CRP.Range("H89").Value = 14886
$ws_CRP.Range("I89").Value = 4501.75
$ws_CRP.Range("K89").Value = 22508.75
$ws_CRP.Range("M89").Value = -16892.75

$ws_CRP.Range("H132").Value = 3298.3333
$ws_CRP.Range("I132").Value = 2982.5
$ws_CRP.Range("K132").Value = 8947.5
$ws_CRP.Range("M132").Value = -6417.5

$ws_CUL.Range("H15").Value = 455052.72
$ws_CUL.Range("I15").Value = 1250294.5
$ws_CUL.Range("J15").Value = 628.8570999999999
$ws_CUL.Range("K15").Value = 3750883.5
$ws_CUL.Range("L15").Value = 1886.5713
$ws_CUL.Range("M15").Value = -3750743.5
$ws_CUL.Range("N15").Value = -2166.5713

$ws_CUL.Range("H75").Value = 3521
$ws_CUL.Range("I75").Value = 1313
$ws_CUL.Range("J75").Value = 4625
$ws_CUL.Range("K75").Value = 3939
$ws_CUL.Range("L75").Value = 13875
$ws_CUL.Range("M75").Value = -2941
$ws_CUL.Range("N75").Value = -15871

$ws_CUL.Range("H78").Value = 3521
$ws_CUL.Range("I78").Value = 1313
$ws_CUL.Range("J78").Value = 4625
$ws_CUL.Range("K78").Value = 11817
$ws_CUL.Range("L78").Value = 41625
$ws_CUL.Range("M78").Value = -6825
$ws_CUL.Range("N78").Value = -51609

$ws_CUL.Range("H107").Value = 611147.5600000001
$ws_CUL.Range("I107").Value = 625.2273
$ws_CUL.Range("K107").Value = 1875.6819
$ws_CUL.Range("M107").Value = 44.31809999999996

$ws_GSM.Range("H97").Value = 1042.2222
$ws_GSM.Range("I97").Value = 912.1111
$ws_GSM.Range("K97").Value = 912.1111
$ws_GSM.Range("M97").Value = -416.1111

$ws_LTW.Range("H132").Value = 2365.3381
$ws_LTW.Range("I132").Value = 2678.4707
$ws_LTW.Range("J132").Value = 2052.2058
$ws_LTW.Range("K132").Value = 8035.4121
$ws_LTW.Range("L132").Value = 6156.617400000001
$ws_LTW.Range("M132").Value = -5505.4121
$ws_LTW.Range("N132").Value = -11216.6174

$ws_WVR.Range("H132").Value = 4573.037
$ws_WVR.Range("I132").Value = 5733.5625
$ws_WVR.Range("J132").Value = 2885
$ws_WVR.Range("K132").Value = 17200.6875
$ws_WVR.Range("L132").Value = 8655
$ws_WVR.Range("M132").Value = -14670.6875
$ws_WVR.Range("N132").Value = -13715
